# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 22:29"

# Update country statistics rows (columns: B=Casos totales, C=Nuevos casos,
# D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8504278
$ws.Range("C4").Value = 47625
$ws.Range("D4").Value = 5530035
$ws.Range("E4").Value = 2748318
$ws.Range("G4").Value = 703
$ws.Range("H4").Value = 225925

# Row 15 - Sudafrica
$ws.Range("B15").Value = 706304
$ws.Range("C15").Value = 1050
$ws.Range("D15").Value = 639568
$ws.Range("E15").Value = 48080
$ws.Range("G15").Value = 164
$ws.Range("H15").Value = 18656

# Row 21 - Alemania
$ws.Range("B21").Value = 380893
$ws.Range("C21").Value = 7162
$ws.Range("E21").Value = 72638
$ws.Range("G21").Value = 56
$ws.Range("H21").Value = 9955

# Row 28 - Israel
$ws.Range("B28").Value = 306162
$ws.Range("C28").Value = 1286
$ws.Range("D28").Value = 282505
$ws.Range("E28").Value = 21379

# Row 31 - Canada
$ws.Range("B31").Value = 203198
$ws.Range("C31").Value = 1761
$ws.Range("D31").Value = 171368
$ws.Range("E31").Value = 22038

# Row 50 - Costa Rica
$ws.Range("B50").Value = 97922
$ws.Range("C50").Value = 847
$ws.Range("D50").Value = 60109
$ws.Range("E50").Value = 36591
$ws.Range("G50").Value = 18
$ws.Range("H50").Value = 1222

# Row 76 - Tunez
$ws.Range("B76").Value = 44450
$ws.Range("C76").Value = 1723
$ws.Range("E76").Value = 38707
$ws.Range("G76").Value = 24
$ws.Range("H76").Value = 711

# Row 104 - Namibia
$ws.Range("B104").Value = 12367
$ws.Range("C104").Value = 41
$ws.Range("D104").Value = 10528
$ws.Range("E104").Value = 1707
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 132

# Row 128 - Republica de Yibuti
$ws.Range("B128").Value = 5499
$ws.Range("C128").Value = 30
$ws.Range("D128").Value = 5384
$ws.Range("E128").Value = 54

# Row 130 - Trinidad yTobago
$ws.Range("B130").Value = 5333
$ws.Range("C130").Value = 35
$ws.Range("D130").Value = 3758
$ws.Range("E130").Value = 1477
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 98

# Row 136 - Ruanda
$ws.Range("B136").Value = 4996
$ws.Range("C136").Value = 4
$ws.Range("E136").Value = 165

# Row 151 - Mali
$ws.Range("B151").Value = 3411
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 2593
$ws.Range("E151").Value = 686

# Row 159 - Sierra Leona
$ws.Range("B159").Value = 2336
$ws.Range("C159").Value = 5
$ws.Range("D159").Value = 1765

# Row 165 - Republica del Chad
$ws.Range("B165").Value = 1399
$ws.Range("C165").Value = 9
$ws.Range("D165").Value = 1199
$ws.Range("E165").Value = 107
